# "Update countries & provincias Spain"
#
# The workbook is a single-sheet ("Pais") COVID-19 country dashboard,
# sorted descending by total cases. This refreshes the data to a later
# snapshot (timestamp 18:22 -> 18:52) which moves Ecuador just ahead of
# Portugal in the ranking (rows 21/22) and Irak just ahead of Armenia
# (rows 68/69), and updates the numeric columns for every country whose
# figures changed between the two snapshots.
#
# Columns: A=Pais B=Casos totales C=Nuevos casos D=Casos activos
#          E=Recuperados F=Casos criticos G=Muertes hoy H=Muertes

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 1 de Mayo de 2020 a las 18:52"

# --- Row 4: Estados Unidos --------------------------------------------
$ws.Range("B4").Value = 1103927
$ws.Range("C4").Value = 8904
$ws.Range("E4").Value = 882948
$ws.Range("G4").Value = 604
$ws.Range("H4").Value = 64460

# --- Row 5: Espana ------------------------------------------------------
$ws.Range("B5").Value = 242988
$ws.Range("C5").Value = 3349
$ws.Range("D5").Value = 142450
$ws.Range("E5").Value = 75714

# --- Row 7: Reino Unido -------------------------------------------------
$ws.Range("B7").Value = 177454
$ws.Range("C7").Value = 6201
$ws.Range("E7").Value = 149600
$ws.Range("G7").Value = 739
$ws.Range("H7").Value = 27510

# --- Row 9: Alemania ------------------------------------------------
$ws.Range("B9").Value = 163542
$ws.Range("C9").Value = 533
$ws.Range("E9").Value = 30010

# --- Row 10: Turquia ------------------------------------------------
$ws.Range("B10").Value = 122392
$ws.Range("C10").Value = 2188
$ws.Range("D10").Value = 53808
$ws.Range("E10").Value = 65326
$ws.Range("F10").Value = 1480
$ws.Range("G10").Value = 84
$ws.Range("H10").Value = 3258

# --- Row 13: Brasil -------------------------------------------------
$ws.Range("B13").Value = 87364
$ws.Range("C13").Value = 1984
$ws.Range("E13").Value = 45412
$ws.Range("G13").Value = 116
$ws.Range("H13").Value = 6017

# --- Row 15: Canada -------------------------------------------------
$ws.Range("B15").Value = 53669
$ws.Range("C15").Value = 433
$ws.Range("D15").Value = 22090
$ws.Range("E15").Value = 28355
$ws.Range("G15").Value = 40
$ws.Range("H15").Value = 3224

# --- Rows 21/22: Ecuador overtakes Portugal -----------------------------
$ws.Range("A21").Value = "Ecuador"
$ws.Range("B21").Value = 26336
$ws.Range("C21").Value = 1402
$ws.Range("D21").Value = 1913
$ws.Range("E21").Value = 23360
$ws.Range("F21").Value = 149
$ws.Range("G21").Value = 163
$ws.Range("H21").Value = 1063

$ws.Range("A22").Value = "Portugal"
$ws.Range("B22").Value = 25351
$ws.Range("C22").Value = 306
$ws.Range("D22").Value = 1647
$ws.Range("E22").Value = 22697
$ws.Range("F22").Value = 154
$ws.Range("G22").Value = 18
$ws.Range("H22").Value = 1007

# --- Row 27: Pakistan -----------------------------------------------
$ws.Range("D27").Value = 4351
$ws.Range("E27").Value = 12854

# --- Row 55: Marruecos ----------------------------------------------
$ws.Range("B55").Value = 4569
$ws.Range("C55").Value = 146
$ws.Range("D55").Value = 1083
$ws.Range("E55").Value = 3315

# --- Row 56: Argentina ----------------------------------------------
$ws.Range("D56").Value = 1292
$ws.Range("E56").Value = 2916
$ws.Range("G56").Value = 2
$ws.Range("H56").Value = 220

# --- Rows 68/69: Irak overtakes Armenia ----------------------------------
$ws.Range("A68").Value = "Irak"
$ws.Range("B68").Value = 2153
$ws.Range("C68").Value = 68
$ws.Range("D68").Value = 1414
$ws.Range("E68").Value = 645
$ws.Range("F68").Value = 0
$ws.Range("H68").Value = 94

$ws.Range("A69").Value = "Armenia"
$ws.Range("B69").Value = 2148
$ws.Range("C69").Value = 82
$ws.Range("D69").Value = 977
$ws.Range("E69").Value = 1138
$ws.Range("F69").Value = 10
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 33

# --- Row 125: Mauricio ----------------------------------------------
$ws.Range("D125").Value = 312
$ws.Range("E125").Value = 10

# --- Row 167: Libia -------------------------------------------------
$ws.Range("B167").Value = 63
$ws.Range("C167").Value = 2
$ws.Range("E167").Value = 42
